# Bosnia Herzegovina Premier Liga base update (30-05-2024 12:21).
# The three most recent fixtures (rows 192-194) were re-synced: their
# identifying id/result/odds data rotates up one row (old row 193 data
# now lives on row 192, old row 194 data now lives on row 193, and old
# row 192 data now lives on row 194), while the rank (col A) and the
# match date (col D) stay put for each row.

$ws = $excel.ActiveWorkbook.ActiveSheet

# --- Row 192 (now FK Sarajevo vs NK Posusje, result D) ---
$ws.Cells.Item(192, 2).Value  = 7952776      # B192 id
$ws.Cells.Item(192, 5).Value  = "FK Sarajevo" # E192 HomeTeam
$ws.Cells.Item(192, 6).Value  = "NK Posusje"  # F192 AwayTeam
$ws.Cells.Item(192, 7).Value  = 1             # G192 FTHG
$ws.Cells.Item(192, 8).Value  = 1             # H192 FTAG
$ws.Cells.Item(192, 9).Value  = 0             # I192 ht_goals_h
$ws.Cells.Item(192, 10).Value = 0             # J192 ht_goals_a
$ws.Cells.Item(192, 11).Value = "D"           # K192 FTR
$ws.Cells.Item(192, 12).Value = 1.571         # L192 oddH_op
$ws.Cells.Item(192, 13).Value = 3.4           # M192 oddD_op
$ws.Cells.Item(192, 14).Value = 5.5           # N192 oddA_op
$ws.Cells.Item(192, 15).Value = 1.363         # O192 oddH
$ws.Cells.Item(192, 16).Value = 3.9           # P192 oddD
$ws.Cells.Item(192, 17).Value = 8             # Q192 oddA
$ws.Cells.Item(192, 18).Value = -1.25         # R192 Ah
$ws.Cells.Item(192, 19).Value = 1.85          # S192 oddAHH
$ws.Cells.Item(192, 20).Value = 1.95          # T192 oddAHA
$ws.Cells.Item(192, 21).Value = 2.75          # U192 AhOU
$ws.Cells.Item(192, 22).Value = 1.925         # V192 oddAHOver
$ws.Cells.Item(192, 23).Value = 1.875         # W192 oddAHUnder
$ws.Cells.Item(192, 24).Value = -1            # X192 PLH
$ws.Cells.Item(192, 25).Value = 2.9           # Y192 PLD
$ws.Cells.Item(192, 26).Value = -1            # Z192 PLA
$ws.Cells.Item(192, 27).Value = -1            # AA192 PL_Ahh
$ws.Cells.Item(192, 28).Value = 0.95          # AB192 PL_Aha
$ws.Cells.Item(192, 29).Value = -1            # AC192 PL_AhOver
$ws.Cells.Item(192, 30).Value = 0.875         # AD192 PL_AhUnder

# --- Row 193 (now Zrinjski Mostar vs FK Tuzla City, result H) ---
$ws.Cells.Item(193, 2).Value  = 7952779        # B193 id
$ws.Cells.Item(193, 5).Value  = "Zrinjski Mostar" # E193 HomeTeam
$ws.Cells.Item(193, 6).Value  = "FK Tuzla City"   # F193 AwayTeam
$ws.Cells.Item(193, 7).Value  = 4              # G193 FTHG
$ws.Cells.Item(193, 8).Value  = 0              # H193 FTAG
$ws.Cells.Item(193, 9).Value  = 2              # I193 ht_goals_h
$ws.Cells.Item(193, 10).Value = 0              # J193 ht_goals_a
$ws.Cells.Item(193, 11).Value = "H"            # K193 FTR
$ws.Cells.Item(193, 12).Value = 1.25           # L193 oddH_op
$ws.Cells.Item(193, 13).Value = 5.75           # M193 oddD_op
$ws.Cells.Item(193, 14).Value = 7              # N193 oddA_op
$ws.Cells.Item(193, 15).Value = 1.055          # O193 oddH
$ws.Cells.Item(193, 16).Value = 13             # P193 oddD
$ws.Cells.Item(193, 17).Value = 17             # Q193 oddA
$ws.Cells.Item(193, 18).Value = -3.5           # R193 Ah
$ws.Cells.Item(193, 19).Value = 1.975          # S193 oddAHH
$ws.Cells.Item(193, 20).Value = 1.825          # T193 oddAHA
$ws.Cells.Item(193, 21).Value = 4.75           # U193 AhOU
$ws.Cells.Item(193, 22).Value = 1.825          # V193 oddAHOver
$ws.Cells.Item(193, 23).Value = 1.975          # W193 oddAHUnder
$ws.Cells.Item(193, 24).Value = 0.05499999999999994 # X193 PLH
$ws.Cells.Item(193, 25).Value = -1             # Y193 PLD
$ws.Cells.Item(193, 26).Value = -1             # Z193 PLA
$ws.Cells.Item(193, 27).Value = 0.9750000000000001 # AA193 PL_Ahh
$ws.Cells.Item(193, 28).Value = -1             # AB193 PL_Aha
$ws.Cells.Item(193, 29).Value = -1             # AC193 PL_AhOver
$ws.Cells.Item(193, 30).Value = 0.9750000000000001 # AD193 PL_AhUnder

# --- Row 194 (now Velez Mostar vs GOSK Gabela, result D) ---
$ws.Cells.Item(194, 2).Value  = 7952780       # B194 id
$ws.Cells.Item(194, 5).Value  = "Velez Mostar" # E194 HomeTeam
$ws.Cells.Item(194, 6).Value  = "GOSK Gabela"  # F194 AwayTeam
$ws.Cells.Item(194, 7).Value  = 3             # G194 FTHG
$ws.Cells.Item(194, 8).Value  = 3             # H194 FTAG
$ws.Cells.Item(194, 9).Value  = 1             # I194 ht_goals_h
$ws.Cells.Item(194, 10).Value = 1             # J194 ht_goals_a
$ws.Cells.Item(194, 11).Value = "D"           # K194 FTR
$ws.Cells.Item(194, 12).Value = 1.4           # L194 oddH_op
$ws.Cells.Item(194, 13).Value = 4             # M194 oddD_op
$ws.Cells.Item(194, 14).Value = 7             # N194 oddA_op
$ws.Cells.Item(194, 15).Value = 1.363         # O194 oddH
$ws.Cells.Item(194, 16).Value = 4.2           # P194 oddD
$ws.Cells.Item(194, 17).Value = 8             # Q194 oddA
$ws.Cells.Item(194, 18).Value = -1.5          # R194 Ah
$ws.Cells.Item(194, 19).Value = 2             # S194 oddAHH
$ws.Cells.Item(194, 20).Value = 1.8           # T194 oddAHA
$ws.Cells.Item(194, 21).Value = 2.75          # U194 AhOU
$ws.Cells.Item(194, 22).Value = 1.825         # V194 oddAHOver
$ws.Cells.Item(194, 23).Value = 1.975         # W194 oddAHUnder
$ws.Cells.Item(194, 24).Value = -1            # X194 PLH
$ws.Cells.Item(194, 25).Value = 3.2           # Y194 PLD
$ws.Cells.Item(194, 26).Value = -1            # Z194 PLA
$ws.Cells.Item(194, 27).Value = -1            # AA194 PL_Ahh
$ws.Cells.Item(194, 28).Value = 0.8           # AB194 PL_Aha
$ws.Cells.Item(194, 29).Value = 0.825         # AC194 PL_AhOver
$ws.Cells.Item(194, 30).Value = -1            # AD194 PL_AhUnder
